$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - values are stored as text, so a leading
# apostrophe is used to force Excel to keep them as text instead of
# auto-converting the numeric-looking strings to numbers.
$ws.Range("D2").Value = "'246.59"
$ws.Range("D4").Value = "'5.462"
$ws.Range("D5").Value = "'0.05664"
$ws.Range("D7").Value = "'0.8008"
$ws.Range("D8").Value = "'1.038"
$ws.Range("D9").Value = "'0.1449"
$ws.Range("D10").Value = "'0.07246"
$ws.Range("D12").Value = "'0.02939"
$ws.Range("D13").Value = "'0.09287"
$ws.Range("D14").Value = "'0.001651"
$ws.Range("D15").Value = "'3.209"
$ws.Range("D16").Value = "'0.04716"
$ws.Range("D17").Value = "'0.0005890"
$ws.Range("D18").Value = "'0.006363"
$ws.Range("D20").Value = "'0.001046"
$ws.Range("D23").Value = "'3.823"
$ws.Range("D24").Value = "'6.420"
$ws.Range("D25").Value = "'2.090"
$ws.Range("D27").Value = "'0.1317"
$ws.Range("D40").Value = "'0.04085"
$ws.Range("D41").Value = "'0.006909"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("D44").Value = "'0.008947"
$ws.Range("D47").Value = "'0.7852"
$ws.Range("D48").Value = "'0.01178"

# Volume(1h) (column E) text updates
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
